# Updated symbol list on Thu Jan 19 22:37:01 UTC 2023 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 with the
# latest scraped values. Values are prefixed with a leading apostrophe so
# Excel stores them as literal text (matching the source data, which keeps
# trailing zeros / exact percentage strings instead of re-formatted numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'293.98"
$ws.Range("E2").Value = "'1.51%"
$ws.Range("D3").Value = "'31.09"
$ws.Range("E3").Value = "'0.91%"
$ws.Range("D4").Value = "'4.940"
$ws.Range("E4").Value = "'1.24%"
$ws.Range("D5").Value = "'0.07337"
$ws.Range("E5").Value = "'2.60%"
$ws.Range("D6").Value = "'2.283"
$ws.Range("E6").Value = "'23.27%"
$ws.Range("D7").Value = "'7.694"
$ws.Range("E7").Value = "'0.73%"
$ws.Range("D8").Value = "'3.769"
$ws.Range("E8").Value = "'1.31%"
$ws.Range("D9").Value = "'0.9124"
$ws.Range("E9").Value = "'1.83%"
$ws.Range("D10").Value = "'0.1687"
$ws.Range("E10").Value = "'2.58%"
$ws.Range("E11").Value = "'8.47%"
$ws.Range("D12").Value = "'0.08285"
$ws.Range("E12").Value = "'1.64%"
$ws.Range("D13").Value = "'0.03104"
$ws.Range("E13").Value = "'3.59%"
$ws.Range("E14").Value = "'0.60%"
$ws.Range("D15").Value = "'0.001521"
$ws.Range("E15").Value = "'1.59%"
$ws.Range("D16").Value = "'0.005751"
$ws.Range("E16").Value = "'-1.40%"
$ws.Range("D17").Value = "'3.481"
$ws.Range("E17").Value = "'0.65%"
$ws.Range("E18").Value = "'-1.34%"
$ws.Range("D19").Value = "'0.3328"
$ws.Range("D20").Value = "'0.1304"
$ws.Range("E20").Value = "'0.89%"
$ws.Range("D21").Value = "'3.969"
$ws.Range("E21").Value = "'-6.97%"
$ws.Range("D22").Value = "'0.2103"
$ws.Range("E22").Value = "'5.03%"
$ws.Range("D23").Value = "'0.04551"
$ws.Range("E23").Value = "'1.74%"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("E24").Value = "'-0.15%"
$ws.Range("D25").Value = "'0.004342"
$ws.Range("E25").Value = "'-6.77%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'3.83%"
$ws.Range("D27").Value = "'0.0003399"
$ws.Range("E27").Value = "'-95.48%"
$ws.Range("D39").Value = "'0.01599"
$ws.Range("E39").Value = "'-2.29%"
$ws.Range("D40").Value = "'0.04444"
$ws.Range("E40").Value = "'2.42%"
$ws.Range("D41").Value = "'0.007338"
$ws.Range("E41").Value = "'-0.61%"
$ws.Range("D42").Value = "'0.008759"
$ws.Range("D43").Value = "'0.1327"
$ws.Range("E43").Value = "'1.64%"
$ws.Range("D44").Value = "'0.002059"
$ws.Range("E44").Value = "'2.60%"
$ws.Range("D45").Value = "'0.009201"
$ws.Range("E45").Value = "'-10.50%"
$ws.Range("D46").Value = "'0.00005946"
$ws.Range("E46").Value = "'1.74%"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("E48").Value = "'1.65%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.04%"
